# Update RF (column I) values for rows 41-79 from the old rate-factor
# (34.2421052631579) to the new 2025 rate-factor (33.09090909090909),
# reflecting the "Update of 2025 data and RF changes" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRF = 33.09090909090909

$ws.Range("I41:I79").Value = $newRF
